$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update to Instabase job data (shifted from old row 6)
$ws.Range("C5").Value = 'SWE, Machine Learning (ML)'
$ws.Range("D5").Value = 'Instabase'
$ws.Range("E5").Value = 'at instabase, we’re passionate about building software to advance the state of the art in computing. we’ve built a fearlessly experimental, customer-obsessed team who are making discoveries to fundamentally change how people build and consume business applications. today, we’re partnering with the world’s leading companies to transform how they use data and technology. if these challenges excite you, we’d love to hear from you!our engineering teamarchitects the underlying operating system, core services, platform infrastructure, dev toolkits, core algorithms, machine learning models, packaged end-user apps, and app store marketplace. instabase engineers are excited to solve hard problems for complex organizations and are self-starters from day one.what you will do:work on a small team of ml engineers to apply machine learning techniques to challenging real world problems such asimage processing / computer visionocrhandwriting recognitionvisual extraction (checkboxes, signatures, radio button, etc.)object detectiontable detectiondocument/text understandingextraction from documentsfield extraction from document with variable structure (for example: paystubs, invoices, forms, etc.)information retrieval from documents with natural languagefinding relevant clause in legal contractsextracting data (for example: effective date, duration, payment terms, etc.) from legal contractsdocument classificationdocument clusteringworking in this area require knowledge in one or more of the following:ms in computer science, engineering, math, science or related field. phd preferred.5+ years industry / academia experienceproven ttrack record of excellence in applying machine learning techniques for solving difficult real world problemsknowledge of deep learning frameworks such as pytorch and tensorflow, computer vision, nlp, or related areas would be a plus.instabase is an equal opportunity employer and values diversity in all forms. instabase does not discriminate on the basis of race, religion, color, national origin, gender identity, sexual orientation, age, marital status, protected veteran status, disability, or any other unlawful factor. instabase also complies with local laws, including the san francisco fair chance ordinance.'
$ws.Range("F5").Value = 'at instabase, we’re passionate about building software to advance the state of the art in computing. we’ve built a fearlessly experimental, customer-obsessed team who are making discoveries to fundamentally change how people build and consume business applications. today, we’re partnering with the world’s leading companies to transform how they use data and technology. if these challenges excite you, we’d love to hear from you!our engineering teamarchitects the underlying operating system, core services, platform infrastructure, dev toolkits, core algorithms, machine learning models, packaged end-user apps, and app store marketplace. instabase engineers are excited to solve hard problems for complex organizations and are self-starters from day one.what you will do:work on a small team of ml engineers to apply machine learning techniques to challenging real world problems such asimage processing / computer visionocrhandwriting recognitionvisual extraction (checkboxes, signatures, radio button, etc.)object detectiontable detectiondocument/text understandingextraction from documentsfield extraction from document with variable structure (for example: paystubs, invoices, forms, etc.)information retrieval from documents with natural languagefinding relevant clause in legal contractsextracting data (for example: effective date, duration, payment terms, etc.) from legal contractsdocument classificationdocument clusteringworking in this area require knowledge in one or more of the following:ms in computer science, engineering, math, science or related field. phd preferred.5+ years industry / academia experienceproven ttrack record of excellence in applying machine learning techniques for solving difficult real world problemsknowledge of deep learning frameworks such as pytorch and tensorflow, computer vision, nlp, or related areas would be a plus.instabase is an equal opportunity employer and values diversity in all forms. instabase does not discriminate on the basis of race, religion, color, national origin, gender identity, sexual orientation, age, marital status, protected veteran status, disability, or any other unlawful factor. instabase also complies with local laws, including the san francisco fair chance ordinance.'
$ws.Range("G5").Value = 0.1270833333333334
$ws.Range("H5").Value = 0.5035714285714287
$ws.Range("I5").Value = '[''build'', ''solve'', ''apply'', ''learn'', ''’re'', ''advance'', ''obsess'', ''make'', ''change'', ''consume'', ''partner'', ''lead'', ''transform'', ''use'', ''excite'', ''love'', ''hear'', ''teamarchitect'', ''underlie'', ''package'', ''do'', ''work'', ''challenge'', ''visionocrhandwrite'', ''languagefinde'', ''contractsextracte'', ''clusteringworke'', ''require'', ''follow'', ''academia'', ''relate'', ''discriminate'', ''protect'', ''include'']'

# Row 6: update to IDEXX job data (shifted from old row 7)
$ws.Range("C6").Value = 'Global Reference Lab Information Systems Intern'
$ws.Range("D6").Value = 'IDEXX'
$ws.Range("E6").Value = 'join us for a 12-week internship with the global reference lab information systems group.idexx reference laboratoriesis a global network united by a shared commitment to enhancing pet care where the true strength in our name is the people behind it. our reference laboratories make it possible for our customers to discover more with our unrelenting commitment to innovation, personalized support, guidance, and expertise while providing the most complete and advanced menu of diagnostic tests along with technology and tools.this team has software product management responsibility for the software applications used by laboratory technicians in our global network of reference laboratories.laboratory technicians within our reference labs use an ecosystem of applications consisting of a mix internally developed plus third party purchased software. thesoftware ecosystemfacilitates taking in the samples, checking them into the labs with sample management, routing the samples through the lab to complete the diagnostics and then sending the result back out to the customer. this team owns software product management for this ecosystem of applications.as software product management within the global reference lab information systems, we bridge the gap between what the users need, and we turn that into requirements and work with the technical teams to prioritize it and then realize it within the laboratory.in this internship role you would…contribute to putting in a regular, proactive patching program for our applications. to start this initiative, what we need is for you to create a catalogue of our applications, when they were last patched, and how many versions they are behind. work with stakeholders to identify maintenance windows when it is a good time to do the patch, and what frequency, etc.another aspect of this internship is to tour a reference laboratory, observe existing processes, and propose ways for new technology to enhance these processes. in this exploratory part of the internship, you would be looking at opportunities for automation, for example, to use voice recognition technology. tools like microsoft hololens are also available for this part of the internship. proposed technologies can include assisted intelligence and machine learning. we would look to you to come up with proposals for innovation.day-to-day asks are likely to includeworking on a patch program with it, software product managers, and business stakeholdersto catalogue one application, for instance, you would go talk to the product owner that is on the team about what the application does – the functionality and the architecture. then you would go to your contact within it to figure out what the patching strategy should be for it. from there you would figure out what we need to do to do the patch, when to schedule it, how frequently they happen, how out of date is it and what other systems does it impact if we go and patch this and not that, so looking at what is our testing and implementation strategy.there are 25+ points within the ecosystem to categorize. by the time you categorize it, you may realize you missed a data point, then go back and get the data point again, then go back to the business about what is realistic for uptime/downtime/when you are going to do the patch, can you coordinate it with the other systems. there are a lot of pieces to it, so we imagine it being very iterative.ideal skills for this role includeexcellent communication skills, willing to speak up and reach out to different people to get things done.knowledge of computer infrastructure, knowledge of systems architecture and application architecture are all helpful for this internshipbeing detail-oriented and organized, with the ability to pull together all the pieces involved in the projectrequirements for the role includecurrently enrolled in an undergraduate or graduate program majoring in computer science, information systems, network architecture, computer engineering, or other related technology majors.applicants must be at least 18 years of age and must have completed at least one year of college.a track record of student success and potential as demonstrated by gpa, research portfolio, prior work experiences and/or the recommendation of a professor preferred.what’s in it for youyou would get exposure to a corporate enterprise it environment and how systems work together and impact each other, and how you establish good governance and maintenance of those systems! you also would have the opportunity to explore implementing the latest technologies in a production laboratory environment.interns will also be invited to participate in an intern-only slate of programming that includes social/networking events.idexx values a diverse workforce and workplace and strongly encourages women, people of color, lgbt individuals, people with disabilities, members of ethnic minorities, foreign-born residents, and veterans to apply.idexx is an equal opportunity employer. applicants will not be discriminated against because of race, color, creed, sex, sexual orientation, gender identity or expression, age, religion, national origin, citizenship status, disability, ancestry, marital status, veteran status, medical condition, or any protected category prohibited by local, state, or federal laws.'
$ws.Range("F6").Value = 'join us for a 12-week internship with the global reference lab information systems group.idexx reference laboratoriesis a global network united by a shared commitment to enhancing pet care where the true strength in our name is the people behind it. our reference laboratories make it possible for our customers to discover more with our unrelenting commitment to innovation, personalized support, guidance, and expertise while providing the most complete and advanced menu of diagnostic tests along with technology and tools.this team has software product management responsibility for the software applications used by laboratory technicians in our global network of reference laboratories.laboratory technicians within our reference labs use an ecosystem of applications consisting of a mix internally developed plus third party purchased software. thesoftware ecosystemfacilitates taking in the samples, checking them into the labs with sample management, routing the samples through the lab to complete the diagnostics and then sending the result back out to the customer. this team owns software product management for this ecosystem of applications.as software product management within the global reference lab information systems, we bridge the gap between what the users need, and we turn that into requirements and work with the technical teams to prioritize it and then realize it within the laboratory.in this internship role you would…contribute to putting in a regular, proactive patching program for our applications. to start this initiative, what we need is for you to create a catalogue of our applications, when they were last patched, and how many versions they are behind. work with stakeholders to identify maintenance windows when it is a good time to do the patch, and what frequency, etc.another aspect of this internship is to tour a reference laboratory, observe existing processes, and propose ways for new technology to enhance these processes. in this exploratory part of the internship, you would be looking at opportunities for automation, for example, to use voice recognition technology. tools like microsoft hololens are also available for this part of the internship. proposed technologies can include assisted intelligence and machine learning. we would look to you to come up with proposals for innovation.day-to-day asks are likely to includeworking on a patch program with it, software product managers, and business stakeholdersto catalogue one application, for instance, you would go talk to the product owner that is on the team about what the application does – the functionality and the architecture. then you would go to your contact within it to figure out what the patching strategy should be for it. from there you would figure out what we need to do to do the patch, when to schedule it, how frequently they happen, how out of date is it and what other systems does it impact if we go and patch this and not that, so looking at what is our testing and implementation strategy.there are 25+ points within the ecosystem to categorize. by the time you categorize it, you may realize you missed a data point, then go back and get the data point again, then go back to the business about what is realistic for uptime/downtime/when you are going to do the patch, can you coordinate it with the other systems. there are a lot of pieces to it, so we imagine it being very iterative.ideal skills for this role includeexcellent communication skills, willing to speak up and reach out to different people to get things done.knowledge of computer infrastructure, knowledge of systems architecture and application architecture are all helpful for this internshipbeing detail-oriented and organized, with the ability to pull together all the pieces involved in the projectrequirements for the role includecurrently enrolled in an undergraduate or graduate program majoring in computer science, information systems, network architecture, computer engineering, or other related technology majors.applicants must be at least 18 years of age and must have completed at least one year of college.a track record of student success and potential as demonstrated by gpa, research portfolio, prior work experiences and/or the recommendation of a professor preferred.what’s in it for youyou would get exposure to a corporate enterprise it environment and how systems work together and impact each other, and how you establish good governance and maintenance of those systems! you also would have the opportunity to explore implementing the latest technologies in a production laboratory environment.interns will also be invited to participate in an intern-only slate of programming that includes social/networking events.idexx values a diverse workforce and workplace and strongly encourages women, people of color, lgbt individuals, people with disabilities, members of ethnic minorities, foreign-born residents, and veterans to apply.idexx is an equal opportunity employer. applicants will not be discriminated against because of race, color, creed, sex, sexual orientation, gender identity or expression, age, religion, national origin, citizenship status, disability, ancestry, marital status, veteran status, medical condition, or any protected category prohibited by local, state, or federal laws.'
$ws.Range("G6").Value = 0.1042272727272727
$ws.Range("H6").Value = 0.364962703962704
$ws.Range("I6").Value = '[''go'', ''do'', ''use'', ''need'', ''patch'', ''look'', ''get'', ''enhance'', ''have'', ''complete'', ''realize'', ''propose'', ''include'', ''figure'', ''impact'', ''categorize'', ''join'', ''unite'', ''share'', ''make'', ''discover'', ''personalize'', ''provide'', ''consist'', ''develop'', ''purchase'', ''ecosystemfacilitate'', ''take'', ''check'', ''route'', ''send'', ''own'', ''bridge'', ''turn'', ''prioritize'', ''contribute'', ''put'', ''start'', ''create'', ''identify'', ''tour'', ''observe'', ''exist'', ''assist'', ''come'', ''includeworke'', ''talk'', ''schedule'', ''happen'', ''miss'', ''coordinate'', ''be'', ''imagine'', ''speak'', ''reach'', ''internshipbee'', ''orient'', ''organize'', ''pull'', ''involve'', ''enrol'', ''major'', ''demonstrate'', ''’s'', ''work'', ''establish'', ''explore'', ''implement'', ''invite'', ''participate'', ''network'', ''encourage'', ''bear'', ''discriminate'', ''protect'', ''prohibit'']'

# Remove old row 7 (data has been consolidated into rows 5-6)
$ws.Rows.Item(7).Delete()